# Week 4 match results entered for the Spring2023Schedule sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Top table (rows 3-10) ---
$ws.Range("E3").Value = "NA"   # Daniel Burcham  - no match played
$ws.Range("E4").Value = "W"    # Leo Hayward     - win
$ws.Range("E5").Value = "W"    # Laura Thompson  - win
$ws.Range("E6").Value = "L"    # Kim Quan        - loss
$ws.Range("E7").Value = "W"    # Scott Berry     - win
$ws.Range("E9").Value = "NA"   # Jason Bohrer    - no match played

# --- Bottom table (rows 15-22) ---
$ws.Range("E15").Value = "NA"  # Jason Bohrer    - no match played
$ws.Range("E17").Value = "W"   # Shelia Lowe     - win
$ws.Range("E18").Value = "W"   # Scott Berry     - win
$ws.Range("E20").Value = "NA"  # Ashley Daniels  - no match played
$ws.Range("E21").Value = "DNP" # Adrian Warden   - did not play
$ws.Range("E22").Value = "L"   # Shakir Donley   - loss

# Move the active selection to reflect where the user was last working.
$ws.Range("F12").Select()
